$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G, shifting old G (d=7) -> H and old H (d=10) -> I
$ws.Range("G:G").Insert()

# Copy the header formatting (bold, border, centered) from F1 onto the new G1 header cell
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New header for inserted column G (d=6)
$ws.Range("G1").Value = "d=6"

# New data values for d=6 column (rows 2-10)
$ws.Range("G2").Value = 43844648418.54645
$ws.Range("G3").Value = 13915906225.55587
$ws.Range("G4").Value = 45017629050.94464
$ws.Range("G5").Value = 42373806532.45045
$ws.Range("G6").Value = 12631326692.04545
$ws.Range("G7").Value = 1496125120.768589
$ws.Range("G8").Value = 9478119032.684235
$ws.Range("G9").Value = 25640696318.56057
$ws.Range("G10").Value = 1378.410655219358
